# cryptos.xlsx – "Updated symbol list on Sat Feb 11 13:59:55 UTC 2023 with GitHub Actions"
#
# The sheet stores every data cell as literal text (numbers, percentages, the
# "--" placeholders, etc. are all plain strings, not real numbers). Excel's
# COM Range.Value setter auto-converts strings that merely *look* numeric
# (e.g. "310.36", "1.15%") into real Number cells, which would change the
# cell type from the source Text type. To keep these as text - exactly like
# the original workbook - values that look numeric are written with a
# leading apostrophe, which is Excel's standard "force text" quote-prefix
# convention and is stripped from the stored value automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, [string]$text) {
    # Leading apostrophe forces Excel to store the value as text even when
    # it looks like a number/percentage ("310.36", "1.15%", "1,206.40%", ...).
    $ws.Range($a1).Value = "'" + $text
}

# Row 2 - BNB
Set-TextValue "D2" "310.36"
Set-TextValue "E2" "1.15%"

# Row 3 - OKB
Set-TextValue "D3" "41.24"
Set-TextValue "E3" "5.65%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.129"
Set-TextValue "E4" "0.55%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07632"
Set-TextValue "E5" "-0.55%"

# Rows 6-17: a new coin (GateToken) moved to the top of this block, shifting
# the remaining coins down by one row; the last coin of the old block (GateToken
# itself, previously row 17) drops out of the visible range below it (LEO slides
# into row 17 instead). Re-assert B (Coin), C (Link), D (Price) and E (Volume).

# Row 6 - GateToken (was FTXToken)
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "4.265"
Set-TextValue "E6" "0.71%"

# Row 7 - FTXToken (was BTSEToken)
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.620"
Set-TextValue "E7" "1.14%"

# Row 8 - BTSEToken (was MXToken)
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.446"
Set-TextValue "E8" "1.48%"

# Row 9 - MXToken (was LiechtensteinCryptoassetsExchange)
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9058"
Set-TextValue "E9" "-1.01%"

# Row 10 - LiechtensteinCryptoassetsExchange (was WazirX)
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1111"
Set-TextValue "E10" "7.45%"

# Row 11 - WazirX (was MandalaExchangeToken)
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1803"
Set-TextValue "E11" "3.24%"

# Row 12 - MandalaExchangeToken (was BitrueCoin)
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09122"
Set-TextValue "E12" "1.25%"

# Row 13 - BitrueCoin (was BitMartToken)
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.04189"
Set-TextValue "E13" "-5.47%"

# Row 14 - BitMartToken (was BitForexToken)
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.1051"
Set-TextValue "E14" "-0.38%"

# Row 15 - BitForexToken (was TigerCash)
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001250"
Set-TextValue "E15" "-0.29%"

# Row 16 - TigerCash (was LEO)
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005661"
Set-TextValue "E16" "-2.43%"

# Row 17 - LEO (was GateToken)
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.342"
Set-TextValue "E17" "-0.44%"

# Row 18 - BitpandaEcosystemToken
Set-TextValue "E18" "0.33%"

# Row 19 - MCDex
Set-TextValue "D19" "6.627"
Set-TextValue "E19" "-5.71%"

# Row 20 - ProBitToken
Set-TextValue "E20" "1.26%"

# Row 21 - ZBToken
Set-TextValue "D21" "0.2734"
Set-TextValue "E21" "-0.90%"

# Row 22 - CoinExToken
Set-TextValue "D22" "0.04031"
Set-TextValue "E22" "-2.53%"

# Row 23 - BitKan
Set-TextValue "E23" "5.05%"

# Row 24 - HotbitToken
Set-TextValue "D24" "0.004065"
Set-TextValue "E24" "-1.01%"

# Row 25 - NitroEx
Set-TextValue "D25" "0.0001299"
Set-TextValue "E25" "-0.12%"

# Row 38 - One
Set-TextValue "D38" "0.02427"
Set-TextValue "E38" "0.96%"

# Row 39 - IDEX
Set-TextValue "D39" "0.05253"
Set-TextValue "E39" "1.07%"

# Row 40 - KickToken
Set-TextValue "D40" "0.007798"
Set-TextValue "E40" "-1.36%"

# Row 41 - BKEXToken
Set-TextValue "D41" "0.1302"
Set-TextValue "E41" "-1.14%"

# Row 42 - Dexo
Set-TextValue "D42" "0.006857"
Set-TextValue "E42" "-3.75%"

# Row 43 - CEJI
Set-TextValue "D43" "0.001949"
Set-TextValue "E43" "-0.09%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.007546"
Set-TextValue "E44" "-10.07%"

# Row 45 - PooCoin
Set-TextValue "D45" "0.3088"
Set-TextValue "E45" "-7.39%"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00006796"
Set-TextValue "E46" "5.66%"

# Row 47 - Kangarootoken
Set-TextValue "E47" "-0.13%"

# Row 48 - BOLO
Set-TextValue "D48" "0.05580"
Set-TextValue "E48" "1,206.40%"

# Row 49 - CoinbaseStockToken
Set-TextValue "E49" "39.90%"

# Row 50 - CryptobidCoin
Set-TextValue "E50" "-0.13%"

# Row 51 - SpecialPowerGold
Set-TextValue "E51" "-0.13%"
